# Saldo_guide.xlsx update:
#  - Re-run date moved from 2024-11-07 to 2024-11-11 (serial 45603 -> 45607)
#    for every data row, reflected in column G and in the sheet/tab name.
#  - A handful of rows had their balance (columns E and H, which always
#    mirror each other in this report) recomputed for the new run date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet/tab to match the new export timestamp.
$ws.Name = "IClientBalance-20241111-101442-"

# 2) Bump the "as of" date in column G for every data row (2..274) from
#    45603 (2024-11-07) to 45607 (2024-11-11).
$firstRow = 2
$lastRow = 274
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45607
}

# 3) Recomputed balances (columns E and H) for the rows whose totals
#    changed between the two runs.
$balanceUpdates = @{
    43  = 21935.75
    49  = 28820.88
    52  = 677.25
    55  = 86.47
    57  = 30.98
    60  = 276.7
    107 = 23356.51
    129 = 2234.33
    218 = 0
    232 = 42744.21
    245 = 26.68
}

foreach ($row in $balanceUpdates.Keys) {
    $newValue = $balanceUpdates[$row]
    $ws.Cells.Item($row, 5).Value = $newValue
    $ws.Cells.Item($row, 8).Value = $newValue
}
